$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New scrape timestamp for this run
$ts = "2025-09-12 18:19:43"

# Rows 2-5: same listings as before, only the "fetched at" timestamp bumps
$ws.Cells.Item(2, 1).Value = $ts
$ws.Cells.Item(3, 1).Value = $ts
$ws.Cells.Item(4, 1).Value = $ts
$ws.Cells.Item(5, 1).Value = $ts

# Row 6
$ws.Cells.Item(6, 1).Value = $ts
$ws.Cells.Item(6, 2).Value = "サブスク型学習サイトの開発"
$ws.Cells.Item(6, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5365024"
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = "◆開発 ◇サイト"

# Row 7
$ws.Cells.Item(7, 1).Value = $ts
$ws.Cells.Item(7, 2).Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Cells.Item(7, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5392235"
$ws.Cells.Item(7, 7).Value = 68
$ws.Cells.Item(7, 8).Value = "◆開発"

# Row 8
$ws.Cells.Item(8, 1).Value = $ts
$ws.Cells.Item(8, 2).Value = "【急募】在庫管理システムの構築!(その後手配管理システムも依頼予定)"
$ws.Cells.Item(8, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5392325"
$ws.Cells.Item(8, 7).Value = 53
$ws.Cells.Item(8, 8).Value = "◇管理"

# Row 9
$ws.Cells.Item(9, 1).Value = $ts
$ws.Cells.Item(9, 2).Value = "Xアカウント(旧Twitter)のスクレイピング(CSV納品)"
$ws.Cells.Item(9, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5392625"
$ws.Cells.Item(9, 7).Value = 40
$ws.Cells.Item(9, 8).Value = "◆スクレイピング"

# Row 10
$ws.Cells.Item(10, 1).Value = $ts
$ws.Cells.Item(10, 2).Value = "【継続案件あり】AWSに精通しているインフラエンジニアを募集します"
$ws.Cells.Item(10, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5392392"
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).ClearContents()

# Row 11
$ws.Cells.Item(11, 1).Value = $ts
$ws.Cells.Item(11, 2).Value = "Vue.jsを使用した「既存ページ修正」+「追加実装」(ピクセルパーフェクト実装)"
$ws.Cells.Item(11, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5392236"
$ws.Cells.Item(11, 7).Value = 18
$ws.Cells.Item(11, 8).ClearContents()

# Row 12
$ws.Cells.Item(12, 1).Value = $ts
$ws.Cells.Item(12, 2).Value = "【急募】A1活用 画像加工とCSV作成のプロを探しています!"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5392360"
$ws.Cells.Item(12, 7).Value = 13
$ws.Cells.Item(12, 8).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = $ts
$ws.Cells.Item(13, 2).Value = "【AWSのプロ募集】事業成長を共に牽引するクラウドインフラの設計・構築パートナー募集中!"
$ws.Cells.Item(13, 4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5392608"
$ws.Cells.Item(13, 7).Value = 10
$ws.Cells.Item(13, 8).ClearContents()

# Row 14
$ws.Cells.Item(14, 1).Value = $ts
$ws.Cells.Item(14, 2).Value = "【急募】エクセルVBAからXLLアドイン作成の依頼"
$ws.Cells.Item(14, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5392307"
$ws.Cells.Item(14, 7).Value = 10
$ws.Cells.Item(14, 8).ClearContents()
